$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3610346189587261
$ws.Range("D2").Value = 0.07094241660363743
$ws.Range("E2").Value = 0.1343833452072793
$ws.Range("F2").Value = 1.605138123829505
$ws.Range("G2").Value = 1.523761136644993
$ws.Range("H2").Value = 1.358990396546005
$ws.Range("K2").Value = 1.607784583014507
$ws.Range("L2").Value = 0.1047187801413045
$ws.Range("M2").Value = 0.5247690533370317
$ws.Range("N2").Value = 1.660301691173793
$ws.Range("C3").Value = 0.3576795385991147
$ws.Range("D3").Value = 0.07144010520130806
$ws.Range("E3").Value = 0.1341124737362378
$ws.Range("F3").Value = 1.594270704478532
$ws.Range("G3").Value = 1.513607927890476
$ws.Range("H3").Value = 1.360833193443241
$ws.Range("K3").Value = 1.482711372977462
$ws.Range("L3").Value = 0.1051046124842223
$ws.Range("M3").Value = 0.4992592473328585
$ws.Range("N3").Value = 1.681417150863711
$ws.Range("C4").Value = 0.3558043240709736
$ws.Range("D4").Value = 0.07175979419485756
$ws.Range("E4").Value = 0.1340088598334006
$ws.Range("F4").Value = 1.588648815512911
$ws.Range("G4").Value = 1.50843570423875
$ws.Range("H4").Value = 1.362697807530949
$ws.Range("K4").Value = 1.406528086871333
$ws.Range("L4").Value = 0.1053833293884558
$ws.Range("M4").Value = 0.4838525054563192
$ws.Range("N4").Value = 1.695020919172014
$ws.Range("C5").Value = 0.3550866115751461
$ws.Range("D5").Value = 0.07189362167367186
$ws.Range("E5").Value = 0.1339823985310531
$ws.Range("F5").Value = 1.586621357945532
$ws.Range("G5").Value = 1.506594102663144
$ws.Range("H5").Value = 1.363641734103169
$ws.Range("K5").Value = 1.375636826409362
$ws.Range("L5").Value = 0.1055074226826029
$ws.Range("M5").Value = 0.4776385979213629
$ws.Range("N5").Value = 1.700725075947268
$ws.Range("C6").Value = 0.3549702409996769
$ws.Range("D6").Value = 0.07191605827042657
$ws.Range("E6").Value = 0.1339789565204015
$ws.Range("F6").Value = 1.586300594700873
$ws.Range("G6").Value = 1.506304353060386
$ws.Range("H6").Value = 1.36380958262572
$ws.Range("K6").Value = 1.370516664082402
$ws.Range("L6").Value = 0.1055286632793049
$ws.Range("M6").Value = 0.4766106768344116
$ws.Range("N6").Value = 1.70168193792539
$ws.Range("C7").Value = 0.3557944566910862
$ws.Range("D7").Value = 0.07176158465435023
$ws.Range("E7").Value = 0.1340084391524456
$ws.Range("F7").Value = 1.588620406527255
$ws.Range("G7").Value = 1.508409791434971
$ws.Range("H7").Value = 1.362709792683688
$ws.Range("K7").Value = 1.406110852637426
$ws.Range("L7").Value = 0.1053849603846793
$ws.Range("M7").Value = 0.4837684415133552
$ws.Range("N7").Value = 1.69509719768723
$ws.Range("C8").Value = 0.3598394129259361
$ws.Range("D8").Value = 0.07111109457699882
$ws.Range("E8").Value = 0.1342769314567462
$ws.Range("F8").Value = 1.601172516338607
$ws.Range("G8").Value = 1.520039365402965
$ws.Range("H8").Value = 1.359473465177757
$ws.Range("K8").Value = 1.564532304567024
$ws.Range("L8").Value = 0.1048431360144022
$ws.Range("M8").Value = 0.5159200469717362
$ws.Range("N8").Value = 1.667449650625279
$ws.Range("C9").Value = 0.369239914524286
$ws.Range("D9").Value = 0.06994720157234369
$ws.Range("E9").Value = 0.1353013821946902
$ws.Range("F9").Value = 1.634160285214293
$ws.Range("G9").Value = 1.55131465235857
$ws.Range("H9").Value = 1.358957890448437
$ws.Range("K9").Value = 1.880066822838785
$ws.Range("L9").Value = 0.1041124944783505
$ws.Range("M9").Value = 0.5810077980499173
$ws.Range("N9").Value = 1.618305956673174
$ws.Range("C10").Value = 0.3770456026659019
$ws.Range("D10").Value = 0.06915989389795563
$ws.Range("E10").Value = 0.1363584900059536
$ws.Range("F10").Value = 1.663556012369725
$ws.Range("G10").Value = 1.579522250880046
$ws.Range("H10").Value = 1.362155191948972
$ws.Range("K10").Value = 2.114903952458235
$ws.Range("L10").Value = 0.1037782441167252
$ws.Range("M10").Value = 0.6300823170338674
$ws.Range("N10").Value = 1.585297910200564
$ws.Range("C11").Value = 0.3807928222815065
$ws.Range("D11").Value = 0.06881639138998885
$ws.Range("E11").Value = 0.1369057211795806
$ws.Range("F11").Value = 1.678061586357501
$ws.Range("G11").Value = 1.593504906866968
$ws.Range("H11").Value = 1.364391100772281
$ws.Range("K11").Value = 2.222403409649814
$ws.Range("L11").Value = 0.1036702289689337
$ws.Range("M11").Value = 0.6526829712687174
$ws.Range("N11").Value = 1.570955541370003
$ws.Range("C12").Value = 0.382240090471953
$ws.Range("D12").Value = 0.06868841904548972
$ws.Range("E12").Value = 0.137122497366132
$ws.Range("F12").Value = 1.683718350014374
$ws.Range("G12").Value = 1.598966381217139
$ws.Range("H12").Value = 1.365350518185693
$ws.Range("K12").Value = 2.26320752800558
$ws.Range("L12").Value = 0.1036356632322821
$ws.Range("M12").Value = 0.6612811236107632
$ws.Range("N12").Value = 1.565621413876443
$ws.Range("C13").Value = 0.3819271370828972
$ws.Range("D13").Value = 0.06871588660055039
$ws.Range("E13").Value = 0.1370753857458666
$ws.Range("F13").Value = 1.68249276530166
$ws.Range("G13").Value = 1.597782730977968
$ws.Range("H13").Value = 1.365138871155722
$ws.Range("K13").Value = 2.254415344093673
$ws.Range("L13").Value = 0.1036428256517112
$ws.Range("M13").Value = 0.6594275881591045
$ws.Range("N13").Value = 1.566765895307181
$ws.Range("C14").Value = 0.3809113229540344
$ws.Range("D14").Value = 0.06880582085830333
$ws.Range("E14").Value = 0.136923363993855
$ws.Range("F14").Value = 1.678523683665958
$ws.Range("G14").Value = 1.593950881764755
$ws.Range("H14").Value = 1.364467771333864
$ws.Range("K14").Value = 2.225758457173242
$ws.Range("L14").Value = 0.1036672581885902
$ws.Range("M14").Value = 0.6533895491075441
$ws.Range("N14").Value = 1.570514753870237
$ws.Range("C15").Value = 0.3802927912202279
$ws.Range("D15").Value = 0.06886118222277204
$ws.Range("E15").Value = 0.1368314905163146
$ws.Range("F15").Value = 1.676113865790654
$ws.Range("G15").Value = 1.591625484040975
$ws.Range("H15").Value = 1.364071394424002
$ws.Range("K15").Value = 2.208217839997928
$ws.Range("L15").Value = 0.1036830492677048
$ws.Range("M15").Value = 0.6496962566444893
$ws.Range("N15").Value = 1.572823678941383
$ws.Range("C16").Value = 0.3768046678952715
$ws.Range("D16").Value = 0.0691826371885611
$ws.Range("E16").Value = 0.1363240631665441
$ws.Range("F16").Value = 1.662630916468288
$ws.Range("G16").Value = 1.578631699673053
$ws.Range("H16").Value = 1.362024825025713
$ws.Range("K16").Value = 2.107892103837855
$ws.Range("L16").Value = 0.1037861897008803
$ws.Range("M16").Value = 0.6286108754471513
$ws.Range("N16").Value = 1.586248781770246
$ws.Range("C17").Value = 0.3747151404467104
$ws.Range("D17").Value = 0.06938358969195946
$ws.Range("E17").Value = 0.1360297737363751
$ws.Range("F17").Value = 1.65465043300469
$ws.Range("G17").Value = 1.570955946199831
$ws.Range("H17").Value = 1.360969718923599
$ws.Range("K17").Value = 2.046517288358359
$ws.Range("L17").Value = 0.1038607452686797
$ws.Range("M17").Value = 0.6157464902479859
$ws.Range("N17").Value = 1.594657211360973
$ws.Range("C18").Value = 0.3735317806277862
$ws.Range("D18").Value = 0.0695005510486908
$ws.Range("E18").Value = 0.1358667505016875
$ws.Range("F18").Value = 1.650166891258976
$ws.Range("G18").Value = 1.5666493185077
$ws.Range("H18").Value = 1.360436380995253
$ws.Range("K18").Value = 2.011279265647886
$ws.Range("L18").Value = 0.1039077723295243
$ws.Range("M18").Value = 0.6083732511212503
$ws.Range("N18").Value = 1.599556864975957
$ws.Range("C19").Value = 0.3731342883041577
$ws.Range("D19").Value = 0.06954038902175519
$ws.Range("E19").Value = 0.1358126257681462
$ws.Range("F19").Value = 1.648667127628926
$ws.Range("G19").Value = 1.565209730703828
$ws.Range("H19").Value = 1.360268420189016
$ws.Range("K19").Value = 1.999359124601824
$ws.Range("L19").Value = 0.103924406594178
$ws.Range("M19").Value = 0.6058812663790718
$ws.Range("N19").Value = 1.601226676955928
$ws.Range("C20").Value = 0.3749356612070471
$ws.Range("D20").Value = 0.06936205526980643
$ws.Range("E20").Value = 0.1360604550139257
$ws.Range("F20").Value = 1.655488928103409
$ws.Range("G20").Value = 1.571761830925198
$ws.Range("H20").Value = 1.361074424027663
$ws.Range("K20").Value = 2.053044207201367
$ws.Range("L20").Value = 0.103852379710407
$ws.Range("M20").Value = 0.6171132333773954
$ws.Range("N20").Value = 1.593755561948182
$ws.Range("C21").Value = 0.3812089244239019
$ws.Range("D21").Value = 0.06877934788074747
$ws.Range("E21").Value = 0.1369677571499111
$ws.Range("F21").Value = 1.679685046381692
$ws.Range("G21").Value = 1.595071860329057
$ws.Range("H21").Value = 1.364661827359328
$ws.Range("K21").Value = 2.234173063456922
$ws.Range("L21").Value = 0.1036599097300943
$ws.Range("M21").Value = 0.6551619878104304
$ws.Range("N21").Value = 1.569410988474417
$ws.Range("C22").Value = 0.3854737065861684
$ws.Range("D22").Value = 0.06841078157916414
$ws.Range("E22").Value = 0.1376164110985165
$ws.Range("F22").Value = 1.696453863022882
$ws.Range("G22").Value = 1.611277564165704
$ws.Range("H22").Value = 1.367663633301333
$ws.Range("K22").Value = 2.35311316092691
$ws.Range("L22").Value = 0.1035710602646276
$ws.Range("M22").Value = 0.6802609045113428
$ws.Range("N22").Value = 1.554065993655852
$ws.Range("C23").Value = 0.3831824156236792
$ws.Range("D23").Value = 0.06860637043138773
$ws.Range("E23").Value = 0.1372651138596055
$ws.Range("F23").Value = 1.687416341349419
$ws.Range("G23").Value = 1.602539055627062
$ws.Range("H23").Value = 1.366001256923568
$ws.Range("K23").Value = 2.289581192222329
$ws.Range("L23").Value = 0.1036150991569578
$ws.Range("M23").Value = 0.6668439149486858
$ws.Range("N23").Value = 1.562204080926678
$ws.Range("C24").Value = 0.3748359079033605
$ws.Range("D24").Value = 0.06937178652356835
$ws.Range("E24").Value = 0.1360465647981179
$ws.Range("F24").Value = 1.655109518886491
$ws.Range("G24").Value = 1.571397159515101
$ws.Range("H24").Value = 1.361026858685278
$ws.Range("K24").Value = 2.050093240243655
$ws.Range("L24").Value = 0.1038561488092213
$ws.Range("M24").Value = 0.6164952582064416
$ws.Range("N24").Value = 1.594162993522605
$ws.Range("C25").Value = 0.3665392256154263
$ws.Range("D25").Value = 0.07025013783667866
$ws.Range("E25").Value = 0.1349708434155303
$ws.Range("F25").Value = 1.624334151200287
$ws.Range("G25").Value = 1.541940298708283
$ws.Range("H25").Value = 1.358470936605045
$ws.Range("K25").Value = 1.794180404928113
$ws.Range("L25").Value = 0.1042745999083223
$ws.Range("M25").Value = 0.5631803952830836
$ws.Range("N25").Value = 1.631057017013766
